$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date-shaped string ("YYYY-MM-DD"). Excel's COM layer
# auto-converts such literals into a real date serial on assignment, so we
# briefly force the cell to Text format, write the literal, then put the
# style back to Normal (no lingering direct formatting on the new cell).
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = "2025-10-03"
$ws.Range("A55").Style = "Normal"

$ws.Range("B55").Value = "15:20:35"
$ws.Range("C55").Value = "1.00 EUR = 1,832.6539"
